$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2..28 from 45420 to 45421
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45421
}

# Delete row 29 entirely (it was removed in the commit)
$ws.Rows.Item(29).Delete()

# Row 28 loses its explicit row height (reverts to default, no customHeight)
$ws.Rows.Item(28).AutoFit()
